$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "analise descritiva - acrescimo das metas"
#
# Insert 5 new columns (F:J) holding the new "meta" metrics
# (meta, meta_avg, meta_std, meta_min, meta_max), placed right after
# "taxa_sucesso" (col E) and before "arrecadado_sucesso" (which, together
# with everything to its right, shifts 5 columns to the right).
# ---------------------------------------------------------------------------

$ws.Range("F1:J1").EntireColumn.Insert()

# --- Header row (row 1) ----------------------------------------------------
$ws.Range("F1").Value = "meta"
$ws.Range("G1").Value = "meta_avg"
$ws.Range("H1").Value = "meta_std"
$ws.Range("I1").Value = "meta_min"
$ws.Range("J1").Value = "meta_max"

# Match the header formatting used by every other header cell: bold font,
# thin border on all sides, centered horizontally and top-aligned vertically.
$ws.Range("F1:J1").Font.Bold = $true
$ws.Range("F1:J1").HorizontalAlignment = -4108
$ws.Range("F1:J1").VerticalAlignment = -4160
$ws.Range("F1:J1").Borders.LineStyle = 1

# --- Data rows ---------------------------------------------------------------
$ws.Range("F2").Value = 13973042.60019265
$ws.Range("G2").Value = 16834.99108456945
$ws.Range("H2").Value = 17015.69760983049
$ws.Range("I2").Value = 31.89582864100442
$ws.Range("J2").Value = 189313.7035611726

$ws.Range("F3").Value = 15599716.7029188
$ws.Range("G3").Value = 11279.62162177787
$ws.Range("H3").Value = 16430.30708090436
$ws.Range("I3").Value = 12.04441558726698
$ws.Range("J3").Value = 198811.9434626772

$ws.Range("F4").Value = 165199.0578149446
$ws.Range("G4").Value = 1205.832538795216
$ws.Range("H4").Value = 2163.288658625352
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 21176.91783511972

# Match the "R$ #,##0.00" currency formatting used by the neighboring
# "arrecadado_*" columns (now shifted to K:S).
$ws.Range("F2:J4").NumberFormat = "R$ #,##0.00"
